$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.297.71'
$ws.Range("E2").Value = '  -1.53%  '
$ws.Range("D3").Value = '1.901.28'
$ws.Range("E3").Value = '  -2.38%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9957'
$ws.Range("E4").Value = '  -0.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.78'
$ws.Range("E5").Value = '  -3.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4615'
$ws.Range("E7").Value = '  -3.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4133'
$ws.Range("E8").Value = '  -0.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.57'
$ws.Range("E9").Value = '  -1.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08013'
$ws.Range("E10").Value = '  -3.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.007'
$ws.Range("E11").Value = '  -3.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.20'
$ws.Range("E12").Value = '  -2.14%  '
$ws.Range("D13").Value = '1.924.56'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.948'
$ws.Range("E14").Value = '  -4.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.113'
$ws.Range("E15").Value = '  -4.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.15'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9947'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001029'
$ws.Range("E18").Value = '  -3.36%  '
$ws.Range("E19").Value = '  -1.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.61'
$ws.Range("E20").Value = '  -2.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").Value = '29.236.64'
$ws.Range("E22").Value = '  -1.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.455'
$ws.Range("E23").Value = '  -2.71%  '
$ws.Range("E24").Value = '  +1.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.197'
$ws.Range("E25").Value = '  -3.80%  '
$ws.Range("D26").Value = '2.129.53'
$ws.Range("E26").Value = '  -1.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.67'
$ws.Range("E27").Value = '  -2.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.69'
$ws.Range("E28").Value = '  -2.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.122'
$ws.Range("E29").Value = '  -3.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.608'
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.46'
$ws.Range("E31").Value = '  -4.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.044'
$ws.Range("E32").Value = '  +1.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09372'
$ws.Range("E33").Value = '  -2.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.426'
$ws.Range("E34").Value = '  -3.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.517'
$ws.Range("E35").Value = '  -4.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.344'
$ws.Range("E36").Value = '  -2.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06078'
$ws.Range("E37").Value = '  -3.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02246'
$ws.Range("E38").Value = '  -3.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.406'
$ws.Range("E39").Value = '  -2.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.176'
$ws.Range("E40").Value = '  -1.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5834'
$ws.Range("E41").Value = '  -4.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9967'
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.19'
$ws.Range("E43").Value = '  -5.14%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1828'
$ws.Range("E44").Value = '  -3.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.243'
$ws.Range("E45").Value = '  -2.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07546'
$ws.Range("E46").Value = '  +2.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.314'
$ws.Range("E47").Value = '  -3.51%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5510'
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.07'
$ws.Range("E49").Value = '  -3.63%  '
$ws.Range("B50").Value = 'PaxosStandard'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.125'
$ws.Range("E50").Value = '  +12.40%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.921'
$ws.Range("E51").Value = '  -4.21%  '
